# Update column G ("K") values for rows 2-15 on the active sheet.
# These values are regenerated save_data (K replacing the old Strike# values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 3
    3  = 8
    4  = 11
    5  = 7
    6  = 1
    7  = 1
    8  = 5
    9  = 5
    10 = 5
    11 = 1
    12 = 1
    13 = 4
    14 = 4
    15 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
